$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Exsisting methods"
$ws.Range("A2").Select()
